# Sierra Leone master data update
# Replace Madagascar (fra / MDG) zone_user_h master data with
# Sierra Leone (eng / SLE) master data, and normalize the casing of the
# resident client service account user id.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "SLE"

$ws.Range("A3").Value = "eng"
$ws.Range("B3").Value = "SLE"
$ws.Range("C3").Value = "service-account-mosip-resident-client"
